$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hybrids")

# --- Step 1: duplicate row 7 (Chorthippus biguttulus / brunneus) into a new
#     row 8, then change Sp1 to Chorthippus mollis -------------------------
$ws.Rows.Item(8).Insert()
for ($col = 1; $col -le 5; $col++) {
    $ws.Cells.Item(8, $col).Value2 = $ws.Cells.Item(7, $col).Value2()
}
$ws.Cells.Item(8, 1).Value2 = "Chorthippus mollis"

# --- Step 2: insert a new "Condition" column before the old "Year" column -
$ws.Columns.Item(4).Insert()
$ws.Cells.Item(1, 4).Value2 = "Condition"

# --- Step 3: fill in the Condition (Laboratory / Field) for every data row
$conditions = @{
    2 = "Laboratory"
    3 = "Laboratory"
    4 = "Field"
    5 = "Field"
    6 = "Field"
    7 = "Field"
    8 = "Field"
    9 = "Laboratory"
    10 = "Laboratory"
    11 = "Field"
    12 = "Field"
    13 = "Field"
    14 = "Field"
    15 = "Field"
    16 = "Field"
    17 = "Field"
    18 = "Laboratory"
    19 = "Field"
    20 = "Field"
    21 = "Field"
    22 = "Field"
    23 = "Field"
    24 = "Field"
    25 = "Laboratory"
    26 = "Laboratory"
    27 = "Field"
    28 = "Field"
    29 = "Field"
    30 = "Laboratory"
    31 = "Laboratory"
    32 = "Laboratory"
    33 = "Field"
}

foreach ($row in $conditions.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $conditions[$row]
}

# --- Step 4: size the new "Condition" column; columns E/F (formerly D/E)
#     keep the widths/bestFit that shifted over automatically with the
#     column insert, so they are intentionally left untouched here. --------
$ws.Columns.Item(4).ColumnWidth = 18.666666666666668
